$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 8 (which holds "extr1"),
# pushing the extr1..extr8 block down to rows 10..17 so we can insert
# "line7" / "line8" rows in their place.
$ws.Range("A8:E9").Insert()

# New "line7" row
$ws.Cells.Item(8, 1).Value = 6
$ws.Cells.Item(8, 1).Font.Bold = $true
$ws.Cells.Item(8, 1).HorizontalAlignment = -4108
$ws.Cells.Item(8, 1).VerticalAlignment = -4160
$ws.Cells.Item(8, 1).Borders.LineStyle = 1
$ws.Cells.Item(8, 2).Value = "line7"
$ws.Cells.Item(8, 3).Value = 14
$ws.Cells.Item(8, 4).Value = 11
$ws.Cells.Item(8, 5).Value = $true

# New "line8" row
$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 1).Font.Bold = $true
$ws.Cells.Item(9, 1).HorizontalAlignment = -4108
$ws.Cells.Item(9, 1).VerticalAlignment = -4160
$ws.Cells.Item(9, 1).Borders.LineStyle = 1
$ws.Cells.Item(9, 2).Value = "line8"
$ws.Cells.Item(9, 3).Value = 16
$ws.Cells.Item(9, 4).Value = 9
$ws.Cells.Item(9, 5).Value = $false

# Updated values for the extr1..extr8 rows, now shifted to rows 10..17.
$extrData = @(
    @{ Row = 10; A = 8;  C = 5;  D = 12; E = $true  },
    @{ Row = 11; A = 9;  C = 5;  D = 9;  E = $true  },
    @{ Row = 12; A = 10; C = 10; D = 11; E = $false },
    @{ Row = 13; A = 11; C = 7;  D = 8;  E = $true  },
    @{ Row = 14; A = 12; C = 9;  D = 11; E = $false },
    @{ Row = 15; A = 13; C = 7;  D = 11; E = $true  },
    @{ Row = 16; A = 14; C = 5;  D = 7;  E = $true  },
    @{ Row = 17; A = 15; C = 8;  D = 5;  E = $false }
)

foreach ($item in $extrData) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
}
